$wb = $excel.ActiveWorkbook

# --- Sheet 1: "SemScores Analysis" ---
$ws1 = $wb.Worksheets.Item("SemScores Analysis")

# Row 24: was boolean FALSE (0) -> category label "Non-numerical"; mean/variance updated
$ws1.Range("B24").Value = "Non-numerical"
$ws1.Range("C24").Value = 0.312853871072685
$ws1.Range("D24").Value = 0.03607403363821524

# Row 25: was boolean TRUE (1) -> category label "Got_supporting_entities"; mean/variance unchanged
$ws1.Range("B25").Value = "Got_supporting_entities"
$ws1.Range("C25").Value = 0.4445277697314547
$ws1.Range("D25").Value = 0.03646076340381461

# New row 26: No_supporting_entities category
$ws1.Range("A26").Value = "got_supporting_ents"
$ws1.Range("B26").Value = "No_supporting_entities"
$ws1.Range("C26").Value = 0.224317392432801
$ws1.Range("D26").Value = 0.02082416114565826

# New row 27: Overall totals
$ws1.Range("A27").Value = "Overall"
$ws1.Range("B27").Value = "ALL"
$ws1.Range("C27").Value = 0.3219779752531655
$ws1.Range("D27").Value = 0.03763840683771095

# --- Sheet 2: "Best Performers" ---
$ws2 = $wb.Worksheets.Item("Best Performers")
$ws2.Range("R2").Value = "Got_supporting_entities"
$ws2.Range("R3").Value = "Got_supporting_entities"
$ws2.Range("R4").Value = "Got_supporting_entities"
$ws2.Range("R5").Value = "Got_supporting_entities"
$ws2.Range("R6").Value = "Got_supporting_entities"
$ws2.Range("R7").Value = "Non-numerical"
$ws2.Range("R8").Value = "Got_supporting_entities"
$ws2.Range("R9").Value = "Got_supporting_entities"
$ws2.Range("R10").Value = "Got_supporting_entities"
$ws2.Range("R11").Value = "Got_supporting_entities"

# --- Sheet 3: "Worst Performers" ---
$ws3 = $wb.Worksheets.Item("Worst Performers")
$ws3.Range("R2").Value = "Non-numerical"
$ws3.Range("R3").Value = "Non-numerical"
$ws3.Range("R4").Value = "Non-numerical"
$ws3.Range("R5").Value = "Non-numerical"
$ws3.Range("R6").Value = "No_supporting_entities"
$ws3.Range("R7").Value = "Non-numerical"
$ws3.Range("R8").Value = "Non-numerical"
$ws3.Range("R9").Value = "Non-numerical"
$ws3.Range("R10").Value = "Non-numerical"
$ws3.Range("R11").Value = "Non-numerical"
